{"js": "// The document's header reads (across 3 paragraphs):\n//   \"V1.0 | NOVEMBER 2024\"\n//   \"Switch Adapted Flap and Wobble Penguin\"\n//   \"Assembly Guide\"\n// The edit renames the third line from \"Assembly Guide\" to \"MAKER Guide\"\n// (displayed in small caps as \"MAKER GUIDE\" because of the paragraph's\n// w:caps run formatting), leaving every other run/paragraph untouched.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\n// Word/Office.js only exposes headers through Section.getHeader(type);\n// \"Primary\" is the default header used on every page of this section\n// (the document does not have distinct first-page/even-page headers).\nconst primaryHeader = section.getHeader(\"Primary\");\n\n// Find the exact run of text we need to retitle and replace it in place so\n// the surrounding run formatting (font, bold, caps, color, size) is kept.\nconst found = primaryHeader.search(\"Assembly Guide\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"MAKER Guide\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document's header reads (across 3 paragraphs):\n#   \"V1.0 | NOVEMBER 2024\"\n#   \"Switch Adapted Flap and Wobble Penguin\"\n#   \"Assembly Guide\"\n# The edit renames the third line from \"Assembly Guide\" to \"MAKER Guide\"\n# (displayed in small caps as \"MAKER GUIDE\" because of the paragraph's\n# caps run formatting), leaving every other run/paragraph untouched.\n\n$d = $word.ActiveDocument\n\n# Use the ordinal indexer (Headers.Item(1)) to reach the section's default\n# header; this only touches the one header story already in the document\n# instead of forcing Word to materialize separate first-page/even-page\n# header parts.\n$section = $d.Sections.Item(1)\n$headerRange = $section.Headers.Item(1).Range\n\n# Scope Find/Replace to that header range only, so nothing in the main\n# document body is touched.\n$find = $headerRange.Find\n$find.ClearFormatting()\n$find.Text = \"Assembly Guide\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"MAKER Guide\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
